$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 153, column A: trim the leading space from " produktivity"
$ws.Range("A153").Value2 = "produktivity"

# 2) Insert a new row at position 507 (shifts rows 507-531 down to 508-532)
#    and populate it with the new "target group" / "QUALITY" / 1 record.
$ws.Rows.Item(507).Insert()
$ws.Range("A507").Value2 = "target group"
$ws.Range("B507").Value2 = "QUALITY"
$ws.Range("C507").Value2 = 1
